$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename sheets
# ---------------------------------------------------------------------------
$wsExercises  = $wb.Worksheets.Item("Exercises")
$wsExercises.Name = "Aufgaben"

$wsThresholds = $wb.Worksheets.Item("Thresholds")
$wsThresholds.Name = "Grenzwerte"

# ---------------------------------------------------------------------------
# 2. Add the new "Module Info" sheet at the end of the workbook
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsModuleInfo = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsModuleInfo.Name = "Module Info"

# ---------------------------------------------------------------------------
# 3. Students sheet - header rename + grading data
# ---------------------------------------------------------------------------
$wsStudents = $wb.Worksheets.Item("Students")

$wsStudents.Cells.Item(1, 4).Value = "Auf 1"
$wsStudents.Cells.Item(1, 5).Value = "Auf 2"
$wsStudents.Cells.Item(1, 6).Value = "Auf 3"
$wsStudents.Cells.Item(1, 7).Value = "Auf 4"

$neText = "                   ne"

# rows whose D:G (Auf1..Auf4) cells must be cleared entirely (no score yet)
$rowsClearDG = @(9, 10, 25, 27)
foreach ($r in $rowsClearDG) {
    foreach ($c in @(4, 5, 6, 7)) {
        $wsStudents.Cells.Item($r, $c).ClearContents()
    }
    $wsStudents.Cells.Item($r, 8).Value = $neText
    $wsStudents.Cells.Item($r, 9).Value = $neText
}

# rows with full numeric results (Auf1..Auf4, Total, Bewertung)
$rowsWithScores = @{
    7  = @(4,  19, 8,    28.5, 59.5, 83)
    8  = @(0,  6,  7.5,  1,    14.5, 21)
    11 = @(0,  1,  0,    0,    1,    1)
    12 = @(7,  11, 16,   17,   51,   75)
    13 = @(2,  12, 10.5, 5,    29.5, 43)
    14 = @(3,  21, 14,   24,   62,   87)
    15 = @(6,  9,  15,   5,    35,   51)
    16 = @(3,  8,  11,   3,    25,   37)
    17 = @(3,  6,  15,   11,   35,   51)
    18 = @(1,  14, 13,   9,    37,   54)
    19 = @(9,  22, 17,   11,   59,   82)
    20 = @(2,  11, 11,   7,    31,   46)
    21 = @(3,  15, 15,   29,   62,   87)
    22 = @(0,  5,  12,   2,    19,   28)
    23 = @(6,  19, 16,   15.5, 56.5, 79)
    24 = @(4,  12, 13,   17,   46,   68)
    26 = @(1,  1,  2.5,  0,    4.5,  7)
}
foreach ($r in $rowsWithScores.Keys) {
    $vals = $rowsWithScores[$r]
    $wsStudents.Cells.Item($r, 4).Value = $vals[0]
    $wsStudents.Cells.Item($r, 5).Value = $vals[1]
    $wsStudents.Cells.Item($r, 6).Value = $vals[2]
    $wsStudents.Cells.Item($r, 7).Value = $vals[3]
    $wsStudents.Cells.Item($r, 8).Value = $vals[4]
    $wsStudents.Cells.Item($r, 9).Value = $vals[5]
}

# rows that keep D:G = 0 (already zero, unchanged) but switch to the
# padded "ne" text in Total/Bewertung
$rowsZeroDGne = @(28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49, 50, 51)
foreach ($r in $rowsZeroDGne) {
    $wsStudents.Cells.Item($r, 4).Value = 0
    $wsStudents.Cells.Item($r, 5).Value = 0
    $wsStudents.Cells.Item($r, 6).Value = 0
    $wsStudents.Cells.Item($r, 7).Value = 0
    $wsStudents.Cells.Item($r, 8).Value = $neText
    $wsStudents.Cells.Item($r, 9).Value = $neText
}

# ---------------------------------------------------------------------------
# 4. Aufgaben sheet (formerly "Exercises")
# ---------------------------------------------------------------------------
$wsExercises.Cells.Item(1, 2).Value = "Anzahl_Punktzahl"

$wsExercises.Cells.Item(2, 1).Value = "Auf 1"
$wsExercises.Cells.Item(3, 1).Value = "Auf 2"
$wsExercises.Cells.Item(4, 1).Value = "Auf 3"
$wsExercises.Cells.Item(5, 1).Value = "Auf 4"

$wsExercises.Cells.Item(3, 2).Value = 25
$wsExercises.Cells.Item(4, 2).Value = 17

# ---------------------------------------------------------------------------
# 5. Grenzwerte sheet (formerly "Thresholds") - append new row 5
# ---------------------------------------------------------------------------
$wsThresholds.Cells.Item(5, 1).Value = 77
$wsThresholds.Cells.Item(5, 2).NumberFormat = "@"
$wsThresholds.Cells.Item(5, 2).Value = "100%"
$wsThresholds.Cells.Item(5, 3).Value = 1

# ---------------------------------------------------------------------------
# 6. Module Info sheet content + column widths
# ---------------------------------------------------------------------------
$wsModuleInfo.Cells.Item(1, 1).Value = "Module Eigenschaft"
$wsModuleInfo.Cells.Item(1, 2).Value = "Module Information"

$wsModuleInfo.Cells.Item(2, 1).Value = "Module Title"
$wsModuleInfo.Cells.Item(2, 2).Value = "Webbasierte Systeme"

$wsModuleInfo.Cells.Item(3, 1).Value = "Module Number"
$wsModuleInfo.Cells.Item(3, 2).Value = "CS1024"

$wsModuleInfo.Cells.Item(4, 1).Value = "Prüfungsdatum"
$wsModuleInfo.Cells.Item(4, 2).NumberFormat = "@"
$wsModuleInfo.Cells.Item(4, 2).Value = "2025-02-27"

$wsModuleInfo.Cells.Item(5, 1).Value = "Prüfer"
$wsModuleInfo.Cells.Item(5, 2).Value = "Herr Müller"

$wsModuleInfo.Cells.Item(6, 1).Value = "Exportdatum"
$wsModuleInfo.Cells.Item(6, 2).NumberFormat = "@"
$wsModuleInfo.Cells.Item(6, 2).Value = "24.2.2025"

$wsModuleInfo.Columns.Item(1).ColumnWidth = 20
$wsModuleInfo.Columns.Item(2).ColumnWidth = 30

Write-Host "edit complete"
